$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates scraped from the authoritative diff: coin prices / 1h volume
# percentages refreshed by the scheduled GitHub Actions run, plus three
# ranking re-orderings (Stacks/USDe/RenderToken and Filecoin/ARBITRUM).
$updates = @(
    @{ Cell = 'D2'; Value = '67.024.28'; ForceText = $false }
    @{ Cell = 'E2'; Value = '  -1.58%  '; ForceText = $false }
    @{ Cell = 'D3'; Value = '2.471.04'; ForceText = $false }
    @{ Cell = 'E3'; Value = '  -2.69%  '; ForceText = $false }
    @{ Cell = 'D5'; Value = '583.73'; ForceText = $true }
    @{ Cell = 'E5'; Value = '  -1.31%  '; ForceText = $false }
    @{ Cell = 'D6'; Value = '169.03'; ForceText = $true }
    @{ Cell = 'E6'; Value = '  -2.56%  '; ForceText = $false }
    @{ Cell = 'E7'; Value = '  +0.07%  '; ForceText = $false }
    @{ Cell = 'E8'; Value = '  -2.04%  '; ForceText = $false }
    @{ Cell = 'D9'; Value = '2.471.26'; ForceText = $false }
    @{ Cell = 'E9'; Value = '  -2.73%  '; ForceText = $false }
    @{ Cell = 'E10'; Value = '  -3.06%  '; ForceText = $false }
    @{ Cell = 'E11'; Value = '  -1.02%  '; ForceText = $false }
    @{ Cell = 'E12'; Value = '  -2.40%  '; ForceText = $false }
    @{ Cell = 'D13'; Value = '0.330'; ForceText = $true }
    @{ Cell = 'E13'; Value = '  -3.33%  '; ForceText = $false }
    @{ Cell = 'D14'; Value = '25.60'; ForceText = $true }
    @{ Cell = 'E14'; Value = '  -3.18%  '; ForceText = $false }
    @{ Cell = 'D15'; Value = '2.924.86'; ForceText = $false }
    @{ Cell = 'E15'; Value = '  -2.03%  '; ForceText = $false }
    @{ Cell = 'D16'; Value = '66.873.57'; ForceText = $false }
    @{ Cell = 'E16'; Value = '  -1.72%  '; ForceText = $false }
    @{ Cell = 'E17'; Value = '  -4.33%  '; ForceText = $false }
    @{ Cell = 'D18'; Value = '2.462.83'; ForceText = $false }
    @{ Cell = 'E18'; Value = '  -2.81%  '; ForceText = $false }
    @{ Cell = 'D19'; Value = '11.10'; ForceText = $true }
    @{ Cell = 'E19'; Value = '  -6.21%  '; ForceText = $false }
    @{ Cell = 'D20'; Value = '7.54'; ForceText = $true }
    @{ Cell = 'E20'; Value = '  -5.92%  '; ForceText = $false }
    @{ Cell = 'D21'; Value = '352.67'; ForceText = $true }
    @{ Cell = 'E21'; Value = '  -5.56%  '; ForceText = $false }
    @{ Cell = 'E22'; Value = '  -2.37%  '; ForceText = $false }
    @{ Cell = 'E23'; Value = '  +0.22%  '; ForceText = $false }
    @{ Cell = 'D24'; Value = '4.25'; ForceText = $true }
    @{ Cell = 'E24'; Value = '  -6.91%  '; ForceText = $false }
    @{ Cell = 'D25'; Value = '68.77'; ForceText = $true }
    @{ Cell = 'E25'; Value = '  -4.48%  '; ForceText = $false }
    @{ Cell = 'D26'; Value = '1.83'; ForceText = $true }
    @{ Cell = 'E26'; Value = '  -4.64%  '; ForceText = $false }
    @{ Cell = 'D27'; Value = '9.19'; ForceText = $true }
    @{ Cell = 'E27'; Value = '  -7.24%  '; ForceText = $false }
    @{ Cell = 'D28'; Value = '0.997'; ForceText = $true }
    @{ Cell = 'E28'; Value = '  -58.53%  '; ForceText = $false }
    @{ Cell = 'D29'; Value = '2.597.74'; ForceText = $false }
    @{ Cell = 'E29'; Value = '  -2.68%  '; ForceText = $false }
    @{ Cell = 'D30'; Value = '0.0₃0905'; ForceText = $false }
    @{ Cell = 'E30'; Value = '  -6.49%  '; ForceText = $false }
    @{ Cell = 'D31'; Value = '511.24'; ForceText = $true }
    @{ Cell = 'E31'; Value = '  -4.80%  '; ForceText = $false }
    @{ Cell = 'D32'; Value = '7.72'; ForceText = $true }
    @{ Cell = 'E32'; Value = '  -7.94%  '; ForceText = $false }
    @{ Cell = 'E33'; Value = '  -6.23%  '; ForceText = $false }
    @{ Cell = 'D34'; Value = '1.77'; ForceText = $true }
    @{ Cell = 'E34'; Value = '  -5.07%  '; ForceText = $false }
    @{ Cell = 'D35'; Value = '1.00'; ForceText = $true }
    @{ Cell = 'E35'; Value = '  -0.09%  '; ForceText = $false }
    @{ Cell = 'D36'; Value = '159.53'; ForceText = $true }
    @{ Cell = 'E36'; Value = '  +0.41%  '; ForceText = $false }
    @{ Cell = 'E37'; Value = '  -10.12%  '; ForceText = $false }
    @{ Cell = 'D38'; Value = '18.66'; ForceText = $true }
    @{ Cell = 'E38'; Value = '  +0.19%  '; ForceText = $false }
    @{ Cell = 'E39'; Value = '  -4.75%  '; ForceText = $false }
    @{ Cell = 'E40'; Value = '  -6.83%  '; ForceText = $false }
    @{ Cell = 'B41'; Value = 'Stacks'; ForceText = $false }
    @{ Cell = 'C41'; Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'; ForceText = $false }
    @{ Cell = 'D41'; Value = '1.70'; ForceText = $true }
    @{ Cell = 'E41'; Value = '  -4.86%  '; ForceText = $false }
    @{ Cell = 'B42'; Value = 'USDe'; ForceText = $false }
    @{ Cell = 'C42'; Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'; ForceText = $false }
    @{ Cell = 'D42'; Value = '1.00'; ForceText = $true }
    @{ Cell = 'E42'; Value = '  -0.20%  '; ForceText = $false }
    @{ Cell = 'B43'; Value = 'RenderToken'; ForceText = $false }
    @{ Cell = 'C43'; Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'; ForceText = $false }
    @{ Cell = 'D43'; Value = '4.84'; ForceText = $true }
    @{ Cell = 'E43'; Value = '  -6.05%  '; ForceText = $false }
    @{ Cell = 'D44'; Value = '0.328'; ForceText = $true }
    @{ Cell = 'E44'; Value = '  -6.63%  '; ForceText = $false }
    @{ Cell = 'E45'; Value = '  -6.84%  '; ForceText = $false }
    @{ Cell = 'D46'; Value = '38.89'; ForceText = $true }
    @{ Cell = 'E46'; Value = '  -1.32%  '; ForceText = $false }
    @{ Cell = 'D47'; Value = '140.86'; ForceText = $true }
    @{ Cell = 'E47'; Value = '  -4.79%  '; ForceText = $false }
    @{ Cell = 'B48'; Value = 'Filecoin'; ForceText = $false }
    @{ Cell = 'C48'; Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; ForceText = $false }
    @{ Cell = 'D48'; Value = '3.45'; ForceText = $true }
    @{ Cell = 'E48'; Value = '  -6.89%  '; ForceText = $false }
    @{ Cell = 'B49'; Value = 'ARBITRUM'; ForceText = $false }
    @{ Cell = 'C49'; Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'; ForceText = $false }
    @{ Cell = 'D49'; Value = '0.514'; ForceText = $true }
    @{ Cell = 'E49'; Value = '  -6.67%  '; ForceText = $false }
    @{ Cell = 'D50'; Value = '0.0₆0256'; ForceText = $false }
    @{ Cell = 'E50'; Value = '  -10.46%  '; ForceText = $false }
    @{ Cell = 'E51'; Value = '  -7.04%  '; ForceText = $false }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    if ($u.ForceText) {
        # Column D normally holds text-formatted price strings (e.g. "583.73").
        # A plain numeric-looking string assigned via .Value would otherwise be
        # auto-coerced to a Number by Excel, dropping formatting like trailing
        # zeros. Force the cell to Text first, then strip the format stamp so
        # no residual style index is left on the cell (matches source which
        # has no "s" attribute on these cells).
        $cell.NumberFormat = "@"
        $cell.Value = $u.Value
        $cell.ClearFormats()
    } else {
        $cell.Value = $u.Value
    }
}
